$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for two new fully-populated rows (music attribution entries) at
# row 9, pushing the old rows 9/10/11 down. The author actually left a gap
# of blank rows (11-14) before continuing the existing "Font/Freeware"
# entries, which end up at rows 15-17.
$ws.Range("9:14").Insert()
$ws.Range("B9:B14").Clear()

# Row 9: The Epic 2 by Rafael Krux
$ws.Range("A9").Value = "the-epic-2-by-rafael-krux.mp3"
$ws.Range("D9").Value = "The Epic 2  by Rafael Krux`nLink: https://filmmusic.io/song/5384-the-epic-2-`nLicense: http://creativecommons.org/licenses/by/4.0/`n Music promoted on https://www.chosic.com/free-music/all/ "
$ws.Range("C9").Value = "Requires Attribution"
$ws.Range("B9").Value = "https://www.chosic.com/download-audio/25862/"

# Row 10: Dragon Slayer by Makai Symphony
$ws.Range("B10").Value = "https://www.chosic.com/download-audio/26014/"
$ws.Range("A10").Value = "makai-symphony-dragon-slayer.mp3"
$ws.Range("C10").Value = "Requires Attribution"
$ws.Range("D10").Value = " Dragon Slayer by Makai Symphony | https://soundcloud.com/makai-symphony`nMusic promoted by https://www.chosic.com/free-music/all/`nCreative Commons Attribution-ShareAlike 3.0 Unported`nhttps://creativecommons.org/licenses/by-sa/3.0/"

# Wrap the long license/notes text and give the two new rows enough height
# to show it.
$ws.Range("C9").WrapText = $true
$ws.Range("D9").WrapText = $true
$ws.Range("D10").WrapText = $true
$ws.Rows.Item(9).RowHeight = 136
$ws.Rows.Item(10).RowHeight = 136

# The row insert doesn't move the worksheet's hyperlink anchors along with
# the cells, so rebuild the hyperlink list pointing at the shifted cells
# (old B10/B11 -> new B16/B17, everything above row 9 is untouched).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B6"), "https://opengameart.org/content/spikes-0") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "https://www.fontspace.com/a-area-kilometer-50-font-f53888") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "https://freesound.org/people/Whiprealgood/sounds/87535/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://freesound.org/people/suntemple/sounds/253172/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "https://opengameart.org/content/simple-explosion-bleeds-game-art") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B7"), "https://opengameart.org/content/various-inventory-24-pixel-icon-set") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B16"), "https://elthen.itch.io/2d-pixel-art-vegetable-monsters-sprite-pack") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B17"), "https://free-game-assets.itch.io/night-city-street-2d-background-tiles") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B8"), "https://opengameart.org/content/energy-icon") | Out-Null

# Match the window scroll position recorded for this edit.
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Range("D11").Select()
